$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and reporting week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/20/2023  Through  11/26/2023"

# --- Weekly crime statistics data updates ---

$ws.Range("C14").Value = "'0"
$ws.Range("C34").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("D15").Value = "'0"
$ws.Range("C34").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "'***.*"
$ws.Range("C34").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 12.5
$ws.Range("I16").Value = 125
$ws.Range("J16").Value = 116
$ws.Range("K16").Value = 7.758620689655
$ws.Range("L16").Value = -10.714285714285
$ws.Range("M16").Value = -60.317460317460
$ws.Range("N16").Value = -87.969201154956
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 175
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 63.157894736842
$ws.Range("I17").Value = 324
$ws.Range("J17").Value = 347
$ws.Range("K17").Value = -6.628242074927
$ws.Range("L17").Value = -11.716621253406
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = -54.172560113154
$ws.Range("C18").Value = "'0"
$ws.Range("C34").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -25
$ws.Range("J18").Value = 116
$ws.Range("K18").Value = -35.344827586206
$ws.Range("L18").Value = -36.974789915966
$ws.Range("M18").Value = -79.619565217391
$ws.Range("N18").Value = -92.588932806324
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -88.888888888888
$ws.Range("F19").Value = 14
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = -58.823529411764
$ws.Range("I19").Value = 342
$ws.Range("J19").Value = 407
$ws.Range("K19").Value = -15.970515970516
$ws.Range("L19").Value = 6.875
$ws.Range("M19").Value = -43.094841930116
$ws.Range("N19").Value = -90.701468189233
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 6.25
$ws.Range("I20").Value = 193
$ws.Range("J20").Value = 212
$ws.Range("K20").Value = -8.962264150943
$ws.Range("L20").Value = 12.865497076023
$ws.Range("M20").Value = -21.862348178137
$ws.Range("N20").Value = -88.210140500916
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -20.833333333333
$ws.Range("F21").Value = 86
$ws.Range("G21").Value = 91
$ws.Range("H21").Value = -5.494505494505
$ws.Range("I21").Value = 1089
$ws.Range("J21").Value = 1241
$ws.Range("K21").Value = -12.248186946011
$ws.Range("L21").Value = -5.221932114882
$ws.Range("M21").Value = -42.532981530343
$ws.Range("N21").Value = -86.690295771205
$ws.Range("F23").Value = "'0"
$ws.Range("C34").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = -100
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 121.428571428571
$ws.Range("F24").Value = 100
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = 12.359550561797
$ws.Range("I24").Value = 1059
$ws.Range("J24").Value = 1185
$ws.Range("K24").Value = -10.632911392405
$ws.Range("L24").Value = 20.477815699658
$ws.Range("M24").Value = 7.186234817813
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -43.75
$ws.Range("F25").Value = 60
$ws.Range("G25").Value = 55
$ws.Range("H25").Value = 9.090909090909
$ws.Range("I25").Value = 617
$ws.Range("J25").Value = 488
$ws.Range("K25").Value = 26.434426229508
$ws.Range("L25").Value = 40.867579908675
$ws.Range("M25").Value = -15.940054495912
$ws.Range("D26").Value = "'0"
$ws.Range("C34").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "'***.*"
$ws.Range("C34").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 33.333333333333
$ws.Range("C36").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("I27").Value = 41
$ws.Range("K27").Value = -21.153846153846
$ws.Range("L27").Value = -8.888888888888
$ws.Range("C36").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 19
$ws.Range("K28").Value = -38.709677419354
$ws.Range("L28").Value = -60.416666666666
$ws.Range("M28").Value = -66.071428571428
$ws.Range("N28").Value = -87.162162162162
$ws.Range("C36").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 14
$ws.Range("K29").Value = -44
$ws.Range("L29").Value = -63.157894736842
$ws.Range("M29").Value = -69.565217391304
$ws.Range("N29").Value = -89.393939393939
